$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder country labels for rows 161-166 so Guinea-Bisau appears
# right after "San Martin (Parte Holandesa)" and before "Islas Caimanes",
# shifting the other countries down one position (Nepal at row167 stays put).
$ws.Range("A161").Value = "Guinea-Bisau"
$ws.Range("A162").Value = "Islas Caimanes"
$ws.Range("A163").Value = "Suazilandia"
$ws.Range("A164").Value = "Benin"
$ws.Range("A165").Value = "Libia"
$ws.Range("A166").Value = "Polinesia Francesa"

# --- Update statistic values (row 23: Irlanda) ---
$ws.Range("B23").Value = 19648
$ws.Range("C23").Value = 386
$ws.Range("E23").Value = 9313
$ws.Range("F23").Value = 146
$ws.Range("G23").Value = 15
$ws.Range("H23").Value = 1102

# --- row 63: Barein ---
$ws.Range("B63").Value = 2723
$ws.Range("C63").Value = 76
$ws.Range("E63").Value = 1497

# --- row 68: Uzbekistan ---
$ws.Range("B68").Value = 1904
$ws.Range("C68").Value = 35
$ws.Range("E68").Value = 1004

# --- row 92: Republica de Chipre ---
$ws.Range("E92").Value = 659
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 15

# --- row 120: Mauricio ---
$ws.Range("E120").Value = 22
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = 10

# --- row 161: now Guinea-Bisau ---
$ws.Range("B161").Value = 73
$ws.Range("C161").Value = 20
$ws.Range("D161").Value = 18
$ws.Range("E161").Value = 54
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 1

# --- row 162: now Islas Caimanes ---
$ws.Range("B162").Value = 70
$ws.Range("C162").Value = 0
$ws.Range("D162").Value = 8
$ws.Range("E162").Value = 61
$ws.Range("F162").Value = 3
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 1

# --- row 163: now Suazilandia ---
$ws.Range("B163").Value = 65
$ws.Range("C163").Value = 6
$ws.Range("D163").Value = 10
$ws.Range("E163").Value = 54
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 0
$ws.Range("H163").Value = 1

# --- row 164: now Benin ---
$ws.Range("B164").Value = 64
$ws.Range("C164").Value = 0
$ws.Range("D164").Value = 33
$ws.Range("E164").Value = 30
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 1

# --- row 165: now Libia ---
$ws.Range("B165").Value = 61
$ws.Range("C165").Value = 0
$ws.Range("D165").Value = 18
$ws.Range("E165").Value = 41
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 2

# --- row 166: now Polinesia Francesa ---
$ws.Range("B166").Value = 57
$ws.Range("C166").Value = 0
$ws.Range("D166").Value = 43
$ws.Range("E166").Value = 14
$ws.Range("F166").Value = 1
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 0
